# Update scripts with new TPM-derived NATMI ligand-receptor values.
# Sending/target cluster labels in columns A and D are re-derived, and the
# numeric ligand/receptor/edge expression + specificity columns (G,H,I,J and
# M..T) are refreshed with the new TPM values for each of the 8 data rows.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 1).Value = "FAPs"
$ws.Cells.Item(2, 4).Value = "ECs"
$ws.Cells.Item(2, 7).Value = 0.5915726666666666
$ws.Cells.Item(2, 8).Value = 1.774718
$ws.Cells.Item(2, 9).Value = 0.95514980205075
$ws.Cells.Item(2, 10).Value = 0.95514980205075
$ws.Cells.Item(2, 13).Value = 0.5550926666666666
$ws.Cells.Item(2, 14).Value = 1.665278
$ws.Cells.Item(2, 15).Value = 0.1208967663154349
$ws.Cells.Item(2, 16).Value = 0.1208967663154349
$ws.Cells.Item(2, 17).Value = 0.328377649067111
$ws.Cells.Item(2, 18).Value = 2.955398841604
$ws.Cells.Item(2, 19).Value = 0.1154745224147634
$ws.Cells.Item(2, 20).Value = 0.1154745224147634
$ws.Cells.Item(3, 1).Value = "FAPs"
$ws.Cells.Item(3, 4).Value = "FAPs"
$ws.Cells.Item(3, 7).Value = 0.5915726666666666
$ws.Cells.Item(3, 8).Value = 1.774718
$ws.Cells.Item(3, 9).Value = 0.95514980205075
$ws.Cells.Item(3, 10).Value = 0.95514980205075
$ws.Cells.Item(3, 13).Value = 3.387303666666666
$ws.Cells.Item(3, 14).Value = 10.161911
$ws.Cells.Item(3, 15).Value = 0.7377399926530269
$ws.Cells.Item(3, 16).Value = 0.7377399926530268
$ws.Cells.Item(3, 17).Value = 2.003836262899778
$ws.Cells.Item(3, 18).Value = 18.034526366098
$ws.Cells.Item(3, 19).Value = 0.7046522079474604
$ws.Cells.Item(3, 20).Value = 0.7046522079474603
$ws.Cells.Item(4, 1).Value = "FAPs"
$ws.Cells.Item(4, 4).Value = "MuSCs"
$ws.Cells.Item(4, 7).Value = 0.5915726666666666
$ws.Cells.Item(4, 8).Value = 1.774718
$ws.Cells.Item(4, 9).Value = 0.95514980205075
$ws.Cells.Item(4, 10).Value = 0.95514980205075
$ws.Cells.Item(4, 13).Value = 0.5311786666666667
$ws.Cells.Item(4, 14).Value = 1.593536
$ws.Cells.Item(4, 15).Value = 0.1156884012202364
$ws.Cells.Item(4, 16).Value = 0.1156884012202364
$ws.Cells.Item(4, 17).Value = 0.3142307803164444
$ws.Cells.Item(4, 18).Value = 2.828077022848
$ws.Cells.Item(4, 19).Value = 0.1104997535250766
$ws.Cells.Item(4, 20).Value = 0.1104997535250765
$ws.Cells.Item(5, 1).Value = "FAPs"
$ws.Cells.Item(5, 4).Value = "Resolving-Mac"
$ws.Cells.Item(5, 7).Value = 0.5915726666666666
$ws.Cells.Item(5, 8).Value = 1.774718
$ws.Cells.Item(5, 9).Value = 0.95514980205075
$ws.Cells.Item(5, 10).Value = 0.95514980205075
$ws.Cells.Item(5, 13).Value = 0.117885
$ws.Cells.Item(5, 14).Value = 0.353655
$ws.Cells.Item(5, 15).Value = 0.02567483981130185
$ws.Cells.Item(5, 16).Value = 0.02567483981130185
$ws.Cells.Item(5, 17).Value = 0.06973754381
$ws.Cells.Item(5, 18).Value = 0.6276378942900001
$ws.Cells.Item(5, 19).Value = 0.02452331816344968
$ws.Cells.Item(5, 20).Value = 0.02452331816344968
$ws.Cells.Item(6, 1).Value = "MuSCs"
$ws.Cells.Item(6, 4).Value = "ECs"
$ws.Cells.Item(6, 7).Value = 0.027778
$ws.Cells.Item(6, 8).Value = 0.08333400000000001
$ws.Cells.Item(6, 9).Value = 0.04485019794925008
$ws.Cells.Item(6, 10).Value = 0.04485019794925008
$ws.Cells.Item(6, 13).Value = 0.5550926666666666
$ws.Cells.Item(6, 14).Value = 1.665278
$ws.Cells.Item(6, 15).Value = 0.1208967663154349
$ws.Cells.Item(6, 16).Value = 0.1208967663154349
$ws.Cells.Item(6, 17).Value = 0.01541936409466666
$ws.Cells.Item(6, 18).Value = 0.138774276852
$ws.Cells.Item(6, 19).Value = 0.005422243900671484
$ws.Cells.Item(6, 20).Value = 0.005422243900671482
$ws.Cells.Item(7, 1).Value = "MuSCs"
$ws.Cells.Item(7, 4).Value = "FAPs"
$ws.Cells.Item(7, 7).Value = 0.027778
$ws.Cells.Item(7, 8).Value = 0.08333400000000001
$ws.Cells.Item(7, 9).Value = 0.04485019794925008
$ws.Cells.Item(7, 10).Value = 0.04485019794925008
$ws.Cells.Item(7, 13).Value = 3.387303666666666
$ws.Cells.Item(7, 14).Value = 10.161911
$ws.Cells.Item(7, 15).Value = 0.7377399926530269
$ws.Cells.Item(7, 16).Value = 0.7377399926530268
$ws.Cells.Item(7, 17).Value = 0.09409252125266666
$ws.Cells.Item(7, 18).Value = 0.8468326912740001
$ws.Cells.Item(7, 19).Value = 0.03308778470556656
$ws.Cells.Item(7, 20).Value = 0.03308778470556655
$ws.Cells.Item(8, 1).Value = "MuSCs"
$ws.Cells.Item(8, 4).Value = "MuSCs"
$ws.Cells.Item(8, 7).Value = 0.027778
$ws.Cells.Item(8, 8).Value = 0.08333400000000001
$ws.Cells.Item(8, 9).Value = 0.04485019794925008
$ws.Cells.Item(8, 10).Value = 0.04485019794925008
$ws.Cells.Item(8, 13).Value = 0.5311786666666667
$ws.Cells.Item(8, 14).Value = 1.593536
$ws.Cells.Item(8, 15).Value = 0.1156884012202364
$ws.Cells.Item(8, 16).Value = 0.1156884012202364
$ws.Cells.Item(8, 17).Value = 0.01475508100266667
$ws.Cells.Item(8, 18).Value = 0.132795729024
$ws.Cells.Item(8, 19).Value = 0.005188647695159868
$ws.Cells.Item(8, 20).Value = 0.005188647695159867
$ws.Cells.Item(9, 1).Value = "MuSCs"
$ws.Cells.Item(9, 4).Value = "Resolving-Mac"
$ws.Cells.Item(9, 7).Value = 0.027778
$ws.Cells.Item(9, 8).Value = 0.08333400000000001
$ws.Cells.Item(9, 9).Value = 0.04485019794925008
$ws.Cells.Item(9, 10).Value = 0.04485019794925008
$ws.Cells.Item(9, 13).Value = 0.117885
$ws.Cells.Item(9, 14).Value = 0.353655
$ws.Cells.Item(9, 15).Value = 0.02567483981130185
$ws.Cells.Item(9, 16).Value = 0.02567483981130185
$ws.Cells.Item(9, 17).Value = 0.00327460953
$ws.Cells.Item(9, 18).Value = 0.02947148577
$ws.Cells.Item(9, 19).Value = 0.001151521647852175
$ws.Cells.Item(9, 20).Value = 0.001151521647852174
